$wb = $excel.ActiveWorkbook

# "Correspond Handoff Datetime" (col E) and "Correspond Handback DateTime" (col H)
# are regenerated for row 2 (the 7520e960... entry) in both the zh-cn and de-de
# handback status sheets, reflecting a fresh report generation run.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 13:28:24"
$wsZhCn.Range("H2").Value = "2016-03-24 13:28:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 13:28:28"
$wsDeDe.Range("H2").Value = "2016-03-24 13:28:56"
